$d = $word.ActiveDocument

# --- 1. Title paragraph: remove leading run before bookmark, merge text into trailing run ---
$titleLead = $d.Content.Duplicate
$titleLead.Find.Execute("Are we ready for Polygenic Risk ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$titleLead.Text = ""

$titleTail = $d.Content.Duplicate
$titleTail.Find.Execute("assessment?", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$titleTail.Text = "Are we ready for Polygenic Risk assessment?"

# --- 2. Background paragraph 1 (FirstParagraph): merge three runs into one ---
$bg1 = "The annual cost of cancer treatment in Australia amounts to billions of dollars worth of burden on the health system. Successful identification and targeted screening of those most at risk offers a far more cost-effective means of reducing cancer burden. Compared to existing risk estimation methods that rely primarily on family-history, Polygenic Risk Assessment offers the potential for a far more precise and personalised mechanism for determining an individual’s cancer risk, even in the absence of a significant family history."
$d.Content.Find.Execute($bg1, $true, $false, $false, $false, $false, $true, 1, $false, $bg1, 2)

# --- 3. Background paragraph 2 (BodyText): merge two runs into one ---
$bg2 = "To determine which malignancies are ready for clinical polygenic risk assessment, we have analysed the NHGRI-EBI catalog of published genome-wide association studies to identify studies reporting malignancy-associated Single Nucleotide Polymorphisms (SNPs) with strong evidence indicating these reported variants are truely associated with the malignancy of interest."
$d.Content.Find.Execute($bg2, $true, $false, $false, $false, $false, $true, 1, $false, $bg2, 2)

# --- 4. Methods paragraph (FirstParagraph): merge three runs into one ---
$mt1 = "Up-to-date publication data for all studies recorded in the GWAS Catalog was extracted from the database, the data were analysed to identify publications reporting SNPs associated with increased risk for various classes of cancer. SNPs for each cancer class were extracted from those identified publications and examined against criteria designed to test the validity of the reported association."
$d.Content.Find.Execute($mt1, $true, $false, $false, $false, $false, $true, 1, $false, $mt1, 2)

# --- 5. Author style: add rsid and rPr (Times New Roman) ---
$authorStyle = $d.Styles("Author")
$authorStyle.Font.Name = "Times New Roman"
